$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Unreported Branches analysis: switch the supplier-specific file name
# references from Humana (Right Source) to RELIANCERX, and populate the
# combined-ins / BO / branch output file names for the new supplier.
$ws.Range("B3").Value = "RELIANCERX"
$ws.Range("B4").Value = "RELIANCERX_SEP23"
$ws.Range("B16").Value = "RELIANCERX Ins"
$ws.Range("B17").Value = "RELIANCERX BO"
$ws.Range("B6").Value = "RELIANCERX BRANCH"
